$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 185, shifting rows 185:268 down to 186:269.
$ws.Rows.Item(185).Insert()

# Populate the new row 185 with the latest weekly price record.
$ws.Cells.Item(185, 1).Value = 6
$ws.Cells.Item(185, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(185, 3).Value = "Metropolitana"
$ws.Cells.Item(185, 4).Value = 44466
$ws.Cells.Item(185, 5).Value = 13
$ws.Cells.Item(185, 6).Value = 100112043
$ws.Cells.Item(185, 7).Value = "Pepino ensalada"
$ws.Cells.Item(185, 8).Value = "Sin especificar"
$ws.Cells.Item(185, 9).Value = "Primera"
$ws.Cells.Item(185, 10).Value = 500
$ws.Cells.Item(185, 11).Value = 14000
$ws.Cells.Item(185, 12).Value = 15000
$ws.Cells.Item(185, 13).Value = 14440
$ws.Cells.Item(185, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(185, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(185, 16).Value = 241
$ws.Cells.Item(185, 17).Value = 60
$ws.Cells.Item(185, 18).Value = "Hortaliza"
